$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The debt-period table in B15:J22 lists periods 2403..2409 (column E) together
# with their overdue amount (column F). The old periods are replaced with the
# new ones: the period list is now given in the opposite order (2409 down to
# 2403), and the single "reduced" amount of 32933 (previously attached to the
# last period, 2409) now travels together with the 2409 row, while every other
# row keeps the regular 52000 amount.

$periods = @("2409", "2408", "2407", "2406", "2405", "2404", "2403")
$amounts = @(32933, 52000, 52000, 52000, 52000, 52000, 52000)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $amounts[$i]
}

$wb.Save()
